$wb = $excel.ActiveWorkbook

# --- Sheet ALC (index 1) ---
$ws = $wb.Worksheets.Item(1)
# row 6
$ws.Range("H6").Value = 729171.4399999999
$ws.Range("J6").Value = 5000
$ws.Range("L6").Value = 15000
$ws.Range("N6").Value = -15224
# row 9
$ws.Range("H9").Value = 58.142857
$ws.Range("I9").Value = 55.4
$ws.Range("J9").Value = 65
$ws.Range("K9").Value = 55.4
$ws.Range("L9").Value = 65
$ws.Range("M9").Value = 113.6
$ws.Range("N9").Value = -403
# row 12
$ws.Range("H12").Value = 153
$ws.Range("I12").Value = 154
$ws.Range("K12").Value = 154
$ws.Range("M12").Value = 16
# row 20
$ws.Range("H20").Value = 4735.8
$ws.Range("I20").Value = 4735.8
$ws.Range("J20").Value = 0
$ws.Range("K20").Value = 4735.8
$ws.Range("L20").Value = 0
$ws.Range("M20").Value = -4505.8
$ws.Range("N20").Value = $null
# row 21
$ws.Range("H21").Value = 4496.75
$ws.Range("I21").Value = 4496.75
$ws.Range("K21").Value = 4496.75
$ws.Range("M21").Value = -4028.75
# row 23
$ws.Range("H23").Value = 4496.75
$ws.Range("I23").Value = 4496.75
$ws.Range("K23").Value = 4496.75
$ws.Range("M23").Value = -4262.75
# row 33
$ws.Range("H33").Value = 174.71428
$ws.Range("I33").Value = 187.16667
$ws.Range("J33").Value = 100
$ws.Range("K33").Value = 187.16667
$ws.Range("L33").Value = 100
$ws.Range("M33").Value = 41.83332999999999
$ws.Range("N33").Value = -558
# row 35
$ws.Range("H35").Value = 4735.8
$ws.Range("I35").Value = 4735.8
$ws.Range("J35").Value = 0
$ws.Range("K35").Value = 4735.8
$ws.Range("L35").Value = 0
$ws.Range("M35").Value = -4356.8
$ws.Range("N35").Value = $null
# row 43
$ws.Range("H43").Value = 5999.5
$ws.Range("J43").Value = 6499.3335
$ws.Range("L43").Value = 6499.3335
$ws.Range("N43").Value = -6637.3335
# row 64
$ws.Range("H64").Value = 3995
$ws.Range("I64").Value = 3995
$ws.Range("K64").Value = 3995
$ws.Range("M64").Value = -3747
# row 67
$ws.Range("H67").Value = 3995
$ws.Range("I67").Value = 3995
$ws.Range("K67").Value = 3995
$ws.Range("M67").Value = -3137
# row 82
$ws.Range("H82").Value = 566.8333
$ws.Range("I82").Value = 566.8333
$ws.Range("K82").Value = 1700.4999
$ws.Range("M82").Value = -1294.4999
# row 85
$ws.Range("H85").Value = 566.8333
$ws.Range("I85").Value = 566.8333
$ws.Range("K85").Value = 1700.4999
$ws.Range("M85").Value = -296.4999
# row 135
$ws.Range("H135").Value = 800
$ws.Range("I135").Value = 800
$ws.Range("J135").Value = 0
$ws.Range("K135").Value = 7200
$ws.Range("L135").Value = 0
$ws.Range("M135").Value = -4665
$ws.Range("N135").Value = $null
# row 137
$ws.Range("H137").Value = 12426.25
$ws.Range("I137").Value = 9186.666999999999
$ws.Range("J137").Value = 15665.833
$ws.Range("K137").Value = 27560.001
$ws.Range("L137").Value = 46997.499
$ws.Range("M137").Value = -25010.001
$ws.Range("N137").Value = -52097.499

# --- Sheet ARM (index 2) ---
$ws = $wb.Worksheets.Item(2)
# row 2
$ws.Range("H2").Value = 4833.3335
$ws.Range("I2").Value = 4000
$ws.Range("J2").Value = 5000
$ws.Range("K2").Value = 4000
$ws.Range("L2").Value = 5000
$ws.Range("M2").Value = -3887
$ws.Range("N2").Value = -5226
# row 92
$ws.Range("H92").Value = 43912.25
$ws.Range("I92").Value = 0
$ws.Range("J92").Value = 43912.25
$ws.Range("K92").Value = 0
$ws.Range("L92").Value = 43912.25
$ws.Range("M92").Value = $null
$ws.Range("N92").Value = -48904.25
# row 98
$ws.Range("H98").Value = 35759.168
$ws.Range("J98").Value = 35759.168
$ws.Range("L98").Value = 35759.168
$ws.Range("N98").Value = -41749.168
# row 116
$ws.Range("H116").Value = 4833.3335
$ws.Range("I116").Value = 4000
$ws.Range("J116").Value = 5000
$ws.Range("K116").Value = 4000
$ws.Range("L116").Value = 5000
$ws.Range("M116").Value = -1706
$ws.Range("N116").Value = -9588

# --- Sheet BSM (index 3) ---
$ws = $wb.Worksheets.Item(3)
# row 3
$ws.Range("H3").Value = 4833.3335
$ws.Range("I3").Value = 4000
$ws.Range("J3").Value = 5000
$ws.Range("K3").Value = 4000
$ws.Range("L3").Value = 5000
$ws.Range("M3").Value = -3886
$ws.Range("N3").Value = -5228
# row 22
$ws.Range("H22").Value = 314.41666
$ws.Range("I22").Value = 287.4
$ws.Range("K22").Value = 287.4
$ws.Range("M22").Value = -114.4
# row 80
$ws.Range("H80").Value = 0
$ws.Range("I80").Value = 0
$ws.Range("J80").Value = 0
$ws.Range("K80").Value = 0
$ws.Range("L80").Value = 0
$ws.Range("M80").Value = $null
$ws.Range("N80").Value = $null
# row 83
$ws.Range("H83").Value = 0
$ws.Range("I83").Value = 0
$ws.Range("J83").Value = 0
$ws.Range("K83").Value = 0
$ws.Range("L83").Value = 0
$ws.Range("M83").Value = $null
$ws.Range("N83").Value = $null

# --- Sheet CRP (index 4) ---
$ws = $wb.Worksheets.Item(4)
# row 31
$ws.Range("H31").Value = 4393
$ws.Range("I31").Value = 3360
$ws.Range("J31").Value = 5839.2
$ws.Range("K31").Value = 3360
$ws.Range("L31").Value = 5839.2
$ws.Range("M31").Value = -3065
$ws.Range("N31").Value = -6429.2
# row 34
$ws.Range("H34").Value = 4393
$ws.Range("I34").Value = 3360
$ws.Range("J34").Value = 5839.2
$ws.Range("K34").Value = 3360
$ws.Range("L34").Value = 5839.2
$ws.Range("M34").Value = -3158
$ws.Range("N34").Value = -6243.2
# row 134
$ws.Range("H134").Value = 2687.3
$ws.Range("I134").Value = 922.7857
$ws.Range("K134").Value = 2768.3571
$ws.Range("M134").Value = -233.3571000000002

# --- Sheet CUL (index 5) ---
$ws = $wb.Worksheets.Item(5)
# row 2
$ws.Range("H2").Value = 274.41666
$ws.Range("J2").Value = 375
$ws.Range("L2").Value = 2250
$ws.Range("N2").Value = -2476
# row 4
$ws.Range("H4").Value = 166668830
$ws.Range("I4").Value = 2599.2
$ws.Range("K4").Value = 7797.599999999999
$ws.Range("M4").Value = -7685.599999999999
# row 11
$ws.Range("H11").Value = 749.6667
$ws.Range("J11").Value = 999.5
$ws.Range("L11").Value = 2998.5
$ws.Range("N11").Value = -3278.5
# row 16
$ws.Range("H16").Value = 250149
$ws.Range("I16").Value = 333433.66
$ws.Range("J16").Value = 295
$ws.Range("K16").Value = 1000300.98
$ws.Range("L16").Value = 885
$ws.Range("M16").Value = -1000127.98
$ws.Range("N16").Value = -1231
# row 34
$ws.Range("H34").Value = 1938.5454
$ws.Range("J34").Value = 2363.7778
$ws.Range("L34").Value = 7091.3334
$ws.Range("N34").Value = -7259.3334

# --- Sheet LTW (index 7) ---
$ws = $wb.Worksheets.Item(7)
# row 22
$ws.Range("H22").Value = 817.75
$ws.Range("I22").Value = 817.75
$ws.Range("K22").Value = 817.75
$ws.Range("M22").Value = -522.75
# row 27
$ws.Range("H27").Value = 817.75
$ws.Range("I27").Value = 817.75
$ws.Range("K27").Value = 817.75
$ws.Range("M27").Value = -710.75
# row 40
$ws.Range("H40").Value = 6833
$ws.Range("I40").Value = 6999.5
$ws.Range("J40").Value = 6500
$ws.Range("K40").Value = 6999.5
$ws.Range("L40").Value = 6500
$ws.Range("M40").Value = -6863.5
$ws.Range("N40").Value = -6772
# row 46
$ws.Range("H46").Value = 6456.6665
$ws.Range("J46").Value = 6353.5713
$ws.Range("L46").Value = 6353.5713
$ws.Range("N46").Value = -6729.5713
# row 136
$ws.Range("H136").Value = 25000
$ws.Range("J136").Value = 25000
$ws.Range("L136").Value = 75000
$ws.Range("N136").Value = -80100

# --- Sheet WVR (index 8) ---
$ws = $wb.Worksheets.Item(8)
# row 28
$ws.Range("H28").Value = 1000
$ws.Range("I28").Value = 1000
$ws.Range("K28").Value = 1000
$ws.Range("M28").Value = -652
# row 104
$ws.Range("H104").Value = 18992
$ws.Range("J104").Value = 18992
$ws.Range("L104").Value = 18992
$ws.Range("N104").Value = -25980
# row 132
$ws.Range("H132").Value = 6493.4116
$ws.Range("I132").Value = 4313
$ws.Range("K132").Value = 12939
$ws.Range("M132").Value = -10409
# row 136
$ws.Range("H136").Value = 7938
$ws.Range("I136").Value = 6680.25
$ws.Range("J136").Value = 18000
$ws.Range("K136").Value = 20040.75
$ws.Range("L136").Value = 54000
$ws.Range("M136").Value = -17490.75
$ws.Range("N136").Value = -59100
